$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete the now-unused last data row (row 10) first so indices below remain valid
$ws.Rows.Item(10).Delete()

# Row 2
$ws.Cells.Item(2,1).Value = "ECs"
$ws.Cells.Item(2,2).Value = "Ccl21b"
$ws.Cells.Item(2,3).Value = "Ccr7"
$ws.Cells.Item(2,4).Value = "Inflammatory-Mac"
$ws.Cells.Item(2,5).Value = 2
$ws.Cells.Item(2,6).Value = 0.6666666666666666
$ws.Cells.Item(2,7).Value = 0.07908133333333334
$ws.Cells.Item(2,8).Value = 0.237244
$ws.Cells.Item(2,9).Value = 0.1380838080781507
$ws.Cells.Item(2,10).Value = 0.1380838080781507
$ws.Cells.Item(2,11).Value = 3
$ws.Cells.Item(2,12).Value = 1
$ws.Cells.Item(2,13).Value = 2.116427666666667
$ws.Cells.Item(2,14).Value = 6.349283
$ws.Cells.Item(2,15).Value = 0.6347244265142036
$ws.Cells.Item(2,16).Value = 0.6347244265142037
$ws.Cells.Item(2,17).Value = 0.1673699217835556
$ws.Cells.Item(2,18).Value = 1.506329296052
$ws.Cells.Item(2,19).Value = 0.08764516589330158
$ws.Cells.Item(2,20).Value = 0.08764516589330158

# Row 3
$ws.Cells.Item(3,1).Value = "ECs"
$ws.Cells.Item(3,2).Value = "Ccl21b"
$ws.Cells.Item(3,3).Value = "Ccr7"
$ws.Cells.Item(3,4).Value = "Resolving-Mac"
$ws.Cells.Item(3,5).Value = 2
$ws.Cells.Item(3,6).Value = 0.6666666666666666
$ws.Cells.Item(3,7).Value = 0.07908133333333334
$ws.Cells.Item(3,8).Value = 0.237244
$ws.Cells.Item(3,9).Value = 0.1380838080781507
$ws.Cells.Item(3,10).Value = 0.1380838080781507
$ws.Cells.Item(3,11).Value = 2
$ws.Cells.Item(3,12).Value = 0.6666666666666666
$ws.Cells.Item(3,13).Value = 1.217976333333333
$ws.Cells.Item(3,14).Value = 3.653929
$ws.Cells.Item(3,15).Value = 0.3652755734857964
$ws.Cells.Item(3,16).Value = 0.3652755734857964
$ws.Cells.Item(3,17).Value = 0.09631919240844444
$ws.Cells.Item(3,18).Value = 0.8668727316760001
$ws.Cells.Item(3,19).Value = 0.05043864218484915
$ws.Cells.Item(3,20).Value = 0.05043864218484915

# Row 4
$ws.Cells.Item(4,1).Value = "FAPs"
$ws.Cells.Item(4,2).Value = "Ccl21b"
$ws.Cells.Item(4,3).Value = "Ccr7"
$ws.Cells.Item(4,4).Value = "Inflammatory-Mac"
$ws.Cells.Item(4,5).Value = 2
$ws.Cells.Item(4,6).Value = 0.6666666666666666
$ws.Cells.Item(4,7).Value = 0.07720866666666666
$ws.Cells.Item(4,8).Value = 0.231626
$ws.Cells.Item(4,9).Value = 0.1348139473702591
$ws.Cells.Item(4,10).Value = 0.134813947370259
$ws.Cells.Item(4,11).Value = 3
$ws.Cells.Item(4,12).Value = 1
$ws.Cells.Item(4,13).Value = 2.116427666666667
$ws.Cells.Item(4,14).Value = 6.349283
$ws.Cells.Item(4,15).Value = 0.6347244265142036
$ws.Cells.Item(4,16).Value = 0.6347244265142037
$ws.Cells.Item(4,17).Value = 0.1634065582397778
$ws.Cells.Item(4,18).Value = 1.470659024158
$ws.Cells.Item(4,19).Value = 0.0855697054307037
$ws.Cells.Item(4,20).Value = 0.0855697054307037

# Row 5
$ws.Cells.Item(5,1).Value = "FAPs"
$ws.Cells.Item(5,2).Value = "Ccl21b"
$ws.Cells.Item(5,3).Value = "Ccr7"
$ws.Cells.Item(5,4).Value = "Resolving-Mac"
$ws.Cells.Item(5,5).Value = 2
$ws.Cells.Item(5,6).Value = 0.6666666666666666
$ws.Cells.Item(5,7).Value = 0.07720866666666666
$ws.Cells.Item(5,8).Value = 0.231626
$ws.Cells.Item(5,9).Value = 0.1348139473702591
$ws.Cells.Item(5,10).Value = 0.134813947370259
$ws.Cells.Item(5,11).Value = 2
$ws.Cells.Item(5,12).Value = 0.6666666666666666
$ws.Cells.Item(5,13).Value = 1.217976333333333
$ws.Cells.Item(5,14).Value = 3.653929
$ws.Cells.Item(5,15).Value = 0.3652755734857964
$ws.Cells.Item(5,16).Value = 0.3652755734857964
$ws.Cells.Item(5,17).Value = 0.09403832872822221
$ws.Cells.Item(5,18).Value = 0.846344958554
$ws.Cells.Item(5,19).Value = 0.04924424193955535
$ws.Cells.Item(5,20).Value = 0.04924424193955534

# Row 6
$ws.Cells.Item(6,1).Value = "Inflammatory-Mac"
$ws.Cells.Item(6,2).Value = "Ccl21b"
$ws.Cells.Item(6,3).Value = "Ccr7"
$ws.Cells.Item(6,4).Value = "Inflammatory-Mac"
$ws.Cells.Item(6,5).Value = 1
$ws.Cells.Item(6,6).Value = 0.3333333333333333
$ws.Cells.Item(6,7).Value = 0.03917266666666667
$ws.Cells.Item(6,8).Value = 0.117518
$ws.Cells.Item(6,9).Value = 0.06839933974190335
$ws.Cells.Item(6,10).Value = 0.06839933974190333
$ws.Cells.Item(6,11).Value = 3
$ws.Cells.Item(6,12).Value = 1
$ws.Cells.Item(6,13).Value = 2.116427666666667
$ws.Cells.Item(6,14).Value = 6.349283
$ws.Cells.Item(6,15).Value = 0.6347244265142036
$ws.Cells.Item(6,16).Value = 0.6347244265142037
$ws.Cells.Item(6,17).Value = 0.08290611551044444
$ws.Cells.Item(6,18).Value = 0.746155039594
$ws.Cells.Item(6,19).Value = 0.04341473169162977
$ws.Cells.Item(6,20).Value = 0.04341473169162977

# Row 7
$ws.Cells.Item(7,1).Value = "Inflammatory-Mac"
$ws.Cells.Item(7,2).Value = "Ccl21b"
$ws.Cells.Item(7,3).Value = "Ccr7"
$ws.Cells.Item(7,4).Value = "Resolving-Mac"
$ws.Cells.Item(7,5).Value = 1
$ws.Cells.Item(7,6).Value = 0.3333333333333333
$ws.Cells.Item(7,7).Value = 0.03917266666666667
$ws.Cells.Item(7,8).Value = 0.117518
$ws.Cells.Item(7,9).Value = 0.06839933974190335
$ws.Cells.Item(7,10).Value = 0.06839933974190333
$ws.Cells.Item(7,11).Value = 2
$ws.Cells.Item(7,12).Value = 0.6666666666666666
$ws.Cells.Item(7,13).Value = 1.217976333333333
$ws.Cells.Item(7,14).Value = 3.653929
$ws.Cells.Item(7,15).Value = 0.3652755734857964
$ws.Cells.Item(7,16).Value = 0.3652755734857964
$ws.Cells.Item(7,17).Value = 0.04771138091355556
$ws.Cells.Item(7,18).Value = 0.429402428222
$ws.Cells.Item(7,19).Value = 0.02498460805027357
$ws.Cells.Item(7,20).Value = 0.02498460805027357

# Row 8
$ws.Cells.Item(8,1).Value = "MuSCs"
$ws.Cells.Item(8,2).Value = "Ccl21b"
$ws.Cells.Item(8,3).Value = "Ccr7"
$ws.Cells.Item(8,4).Value = "Inflammatory-Mac"
$ws.Cells.Item(8,5).Value = 3
$ws.Cells.Item(8,6).Value = 1
$ws.Cells.Item(8,7).Value = 0.3772426666666667
$ws.Cells.Item(8,8).Value = 1.131728
$ws.Cells.Item(8,9).Value = 0.658702904809687
$ws.Cells.Item(8,10).Value = 0.6587029048096869
$ws.Cells.Item(8,11).Value = 3
$ws.Cells.Item(8,12).Value = 1
$ws.Cells.Item(8,13).Value = 2.116427666666667
$ws.Cells.Item(8,14).Value = 6.349283
$ws.Cells.Item(8,15).Value = 0.6347244265142036
$ws.Cells.Item(8,16).Value = 0.6347244265142037
$ws.Cells.Item(8,17).Value = 0.7984068167804445
$ws.Cells.Item(8,18).Value = 7.185661351024001
$ws.Cells.Item(8,19).Value = 0.4180948234985686
$ws.Cells.Item(8,20).Value = 0.4180948234985686

# Row 9
$ws.Cells.Item(9,1).Value = "MuSCs"
$ws.Cells.Item(9,2).Value = "Ccl21b"
$ws.Cells.Item(9,3).Value = "Ccr7"
$ws.Cells.Item(9,4).Value = "Resolving-Mac"
$ws.Cells.Item(9,5).Value = 3
$ws.Cells.Item(9,6).Value = 1
$ws.Cells.Item(9,7).Value = 0.3772426666666667
$ws.Cells.Item(9,8).Value = 1.131728
$ws.Cells.Item(9,9).Value = 0.658702904809687
$ws.Cells.Item(9,10).Value = 0.6587029048096869
$ws.Cells.Item(9,11).Value = 2
$ws.Cells.Item(9,12).Value = 0.6666666666666666
$ws.Cells.Item(9,13).Value = 1.217976333333333
$ws.Cells.Item(9,14).Value = 3.653929
$ws.Cells.Item(9,15).Value = 0.3652755734857964
$ws.Cells.Item(9,16).Value = 0.3652755734857964
$ws.Cells.Item(9,17).Value = 0.4594726399235556
$ws.Cells.Item(9,18).Value = 4.135253759312
$ws.Cells.Item(9,19).Value = 0.2406080813111183
$ws.Cells.Item(9,20).Value = 0.2406080813111183
